$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

$ws.Range('D2').Value = '59.321.12'
$ws.Range('E2').Value = '  -7.05%  '
$ws.Range('D3').Value = '3.298.32'
$ws.Range('E3').Value = '  -4.09%  '
$ws.Range('E4').Value = '  -0.05%  '
Set-TextValue 'D5' '555.41'
$ws.Range('E5').Value = '  -4.61%  '
$ws.Range('E6').Value = '  -2.46%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = '3.296.29'
$ws.Range('E8').Value = '  -4.19%  '
$ws.Range('E9').Value = '  -2.66%  '
Set-TextValue 'D10' '7.31'
$ws.Range('E10').Value = '  -4.33%  '
$ws.Range('E11').Value = '  -6.05%  '
$ws.Range('E12').Value = '  -3.93%  '
$ws.Range('D13').Value = '3.854.60'
$ws.Range('E13').Value = '  -4.18%  '
$ws.Range('E14').Value = '  -0.04%  '
$ws.Range('D15').Value = '3.296.36'
$ws.Range('E15').Value = '  -4.16%  '
$ws.Range('E16').Value = '  -6.37%  '
$ws.Range('E17').Value = '  -4.52%  '
$ws.Range('D18').Value = '59.436.49'
$ws.Range('E18').Value = '  -6.71%  '
Set-TextValue 'D19' '5.59'
$ws.Range('E19').Value = '  -1.34%  '
$ws.Range('E20').Value = '  -1.43%  '
Set-TextValue 'D21' '8.86'
$ws.Range('E21').Value = '  -10.33%  '
Set-TextValue 'D22' '348.71'
$ws.Range('E22').Value = '  -9.01%  '
Set-TextValue 'D23' '0.550'
$ws.Range('E23').Value = '  -2.43%  '
$ws.Range('E24').Value = '  +0.40%  '
$ws.Range('D25').Value = '3.420.18'
$ws.Range('E25').Value = '  -4.33%  '
Set-TextValue 'D26' '68.37'
$ws.Range('E26').Value = '  -7.53%  '
$ws.Range('E27').Value = '  -1.96%  '
Set-TextValue 'D28' '1.00'
$ws.Range('E28').Value = '  +0.50%  '
Set-TextValue 'D29' '7.23'
$ws.Range('E29').Value = '  +2.63%  '
Set-TextValue 'D30' '1.45'
$ws.Range('E30').Value = '  +2.85%  '
$ws.Range('E31').Value = '  -2.32%  '
Set-TextValue 'D32' '7.74'
$ws.Range('E32').Value = '  -2.31%  '
$ws.Range('B33').Value = 'USDe'
$ws.Range('C33').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue 'D33' '1.00'
$ws.Range('E33').Value = '  +0.00%  '
$ws.Range('B34').Value = 'PancakeSwap'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue 'D34' '2.07'
$ws.Range('E34').Value = '  -6.22%  '
$ws.Range('D35').Value = '3.317.26'
$ws.Range('E35').Value = '  -4.28%  '
Set-TextValue 'D36' '22.54'
$ws.Range('E36').Value = '  -1.33%  '
Set-TextValue 'D37' '5.28'
$ws.Range('E37').Value = '  +1.82%  '
Set-TextValue 'D38' '6.73'
$ws.Range('E38').Value = '  -0.39%  '
$ws.Range('E39').Value = '  -1.75%  '
Set-TextValue 'D40' '156.68'
$ws.Range('E40').Value = '  -4.13%  '
Set-TextValue 'D41' '0.0741'
$ws.Range('E41').Value = '  -4.03%  '
$ws.Range('E42').Value = '  -0.18%  '
Set-TextValue 'D43' '40.50'
$ws.Range('E43').Value = '  -1.89%  '
$ws.Range('B44').Value = 'Mantle'
$ws.Range('C44').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue 'D44' '0.738'
$ws.Range('E44').Value = '  -6.73%  '
$ws.Range('E45').Value = '  -1.53%  '
$ws.Range('B46').Value = 'ONDO'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
Set-TextValue 'D46' '1.16'
$ws.Range('E46').Value = '  +4.91%  '
Set-TextValue 'D47' '22.61'
$ws.Range('E47').Value = '  -3.02%  '
$ws.Range('E48').Value = '  -5.29%  '
Set-TextValue 'D49' '6.69'
$ws.Range('E49').Value = '  -0.29%  '
Set-TextValue 'D50' '2.37'
$ws.Range('E50').Value = '  +14.55%  '
Set-TextValue 'D51' '21.63'
$ws.Range('E51').Value = '  +6.30%  '
